$d = $word.ActiveDocument

function Replace-Text($findText, $replaceText) {
    $rng = $d.Content
    $rng.Find.ClearFormatting()
    $found = $rng.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)
    if (-not $found) {
        Write-Output "NOT FOUND: $findText"
    }
    return $found
}

# ---------------------------------------------------------------------------
# 1) Opening sentence: "This vignette using example ... (2018) to generate..."
#    -> "This vignette provides steps for using example ... (submitted) to generate..."
# ---------------------------------------------------------------------------
Replace-Text "This vignette using example" "This vignette provides steps for using example"
Replace-Text "Tracy et al. (2018) to generate" "Tracy et al. (submitted) to generate"

# ---------------------------------------------------------------------------
# 2) "1) Open empty ArcMap project" -> "1) Open a new empty ArcMap project"
# ---------------------------------------------------------------------------
Replace-Text "1) Open empty ArcMap project" "1) Open a new empty ArcMap project"

# ---------------------------------------------------------------------------
# 3) "for setting environments" -> "for setting the Geoprocessing Environments"
# ---------------------------------------------------------------------------
Replace-Text "for setting environments" "for setting the Geoprocessing Environments"

# ---------------------------------------------------------------------------
# 4) Insert new clause before the period ending "...in the above ArcMap project."
#    " so that the output coordinates, processing extent, snap raster and cell
#    size set are set to match the raster " + italic "pop10kmn3"
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.ClearFormatting()
$found = $rng.Find.Execute("in the above ArcMap project", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $ip = $d.Range($rng.End, $rng.End)
    $ip.InsertBefore(" so that the output coordinates, processing extent, snap raster and cell size set are set to match the raster ")
    $ip.Collapse(0)
    $ip.InsertBefore("pop10kmn3")
    $ip.Font.Italic = $true
} else {
    Write-Output "NOT FOUND: in the above ArcMap project"
}

# ---------------------------------------------------------------------------
# 5) "paste in several lines at a time of commands from the"
#    -> "paste in several lines at a time of commands to run from the"
# ---------------------------------------------------------------------------
Replace-Text "commands from the ArcPython script" "commands to run from the ArcPython script"

# ---------------------------------------------------------------------------
# 6) Remove the _GoBack bookmark located between "upweight monarch roost " and
#    "data in low population areas" (text itself is unaffected by this step;
#    handled together with step 9 when the bookmark gets re-added elsewhere).
# ---------------------------------------------------------------------------
try {
    $bm = $d.Bookmarks("_GoBack")
    $bm.Delete()
} catch {
    Write-Output "no _GoBack bookmark found to remove"
}

# ---------------------------------------------------------------------------
# 7) Rename the R script file mentioned in "open the R script
#    KDEMSubset_GridTrainTestEvalAIC_Calib_Function.R" -> "KDEModel_MonarchRoostTSE.R"
#    (keeps italics, since it replaces text within the already-italic run)
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.ClearFormatting()
$found = $rng.Find.Execute("KDEMSubset_GridTrainTestEvalAIC_Calib_Function.R", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $rng.Text = "KDEModel_MonarchRoostTSE.R"
} else {
    Write-Output "NOT FOUND: KDEMSubset_GridTrainTestEvalAIC_Calib_Function.R (first occurrence)"
}

# ---------------------------------------------------------------------------
# 8) Right after that run, the text "(located in the " becomes " (located in the "
#    i.e. an extra leading space is introduced. We simply insert a space before "(located".
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.ClearFormatting()
$found = $rng.Find.Execute("(located in the", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $ip = $d.Range($rng.Start, $rng.Start)
    $ip.InsertBefore(" ")
} else {
    Write-Output "NOT FOUND: (located in the"
}

# ---------------------------------------------------------------------------
# 9) "folder). Make sure the InDirect folders is set to match the desired
#    location on your computer (line 56). Run the R script..."
#    becomes a much longer passage. We break this into ordered pieces.
# ---------------------------------------------------------------------------

# 9a) "folder). Make sure " -> "folder). Make sure FunctDirect matches the location of "
Replace-Text "folder). Make sure the InDirect" "folder). Make sure FunctDirect matches the location of KDEMSubset_GridTrainTestEvalAIC_Calib_Function.R  on your computer. Also, make sure the InDirect"

# 9b) Italicize the newly-inserted "KDEMSubset_GridTrainTestEvalAIC_Calib_Function.R" occurrence
#     (the second occurrence in the document; first occurrence was already renamed in step 7).
$rng = $d.Content
$rng.Find.ClearFormatting()
$found = $rng.Find.Execute("KDEMSubset_GridTrainTestEvalAIC_Calib_Function.R", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $rng.Font.Italic = $true
} else {
    Write-Output "NOT FOUND: KDEMSubset_GridTrainTestEvalAIC_Calib_Function.R (second occurrence)"
}

# 9c) " folders is set to match the desired location"
#     -> " folders in the script are matching the desired location"
Replace-Text " folders is set to match the desired location" " folders in the script are matching the desired location"

# ---------------------------------------------------------------------------
# 10) Re-add the _GoBack bookmark immediately before "Run the R script from
#     lines 1 to 128"
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.ClearFormatting()
$found = $rng.Find.Execute("Run the R script from lines 1 to 128", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $ip = $d.Range($rng.Start, $rng.Start)
    $d.Bookmarks.Add("_GoBack", $ip)
} else {
    Write-Output "NOT FOUND: Run the R script from lines 1 to 128"
}

Write-Output "DONE"
